# The workbook originally has two tabs: "variables" (first/active-ish) and
# "Sheet2" (second, the actually-selected tab). This edit renames "Sheet2"
# to "steps" and drags it in front of "variables", making it the first and
# selected sheet (matching the author's re-upload of the workbook with the
# sheets reordered/renamed).

$wb = $excel.ActiveWorkbook

# Rename the "Sheet2" tab to "steps".
$stepsSheet = $wb.Worksheets.Item("Sheet2")
$stepsSheet.Name = "steps"

# Move it so it becomes the first sheet (before "variables").
$stepsSheet.Move($wb.Worksheets.Item(1))

# Make sure "steps" is the selected/active tab after the reorder. Look it
# up fresh (by name) rather than reusing the pre-move object reference, so
# the activation actually lands on the sheet's new position.
$wb.Worksheets.Item("steps").Activate()
